$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 82, shifting existing rows 82+ down by one.
$ws.Rows.Item(82).Insert()

# Populate the new row 82 with the new translation entry.
$ws.Cells.Item(82, 1).Value = "map_latest_measurements"
$ws.Cells.Item(82, 3).Value = "Latest measurements"
$ws.Cells.Item(82, 4).Value = "Plus récentes mesures"

# Select the new row as the active selection (matches authored selection).
$ws.Rows.Item(82).Select()
